$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "66.821.34"
Set-TextValue "E2" "  -1.44%  "
Set-TextValue "D3" "2.596.70"
Set-TextValue "E3" "  -0.76%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "588.93"
Set-TextValue "E5" "  -2.08%  "
Set-TextValue "D6" "149.55"
Set-TextValue "E6" "  -3.24%  "
Set-TextValue "E7" "  +0.07%  "
Set-TextValue "E8" "  -1.40%  "
Set-TextValue "D9" "2.595.45"
Set-TextValue "E9" "  -0.70%  "
Set-TextValue "E10" "  -2.52%  "
Set-TextValue "E12" "  -1.66%  "
Set-TextValue "E13" "  -3.00%  "
Set-TextValue "D14" "27.17"
Set-TextValue "E14" "  -1.66%  "
Set-TextValue "D15" "3.068.91"
Set-TextValue "E15" "  -0.73%  "
Set-TextValue "E16" "  -5.20%  "
Set-TextValue "D17" "66.778.20"
Set-TextValue "E17" "  -1.36%  "
Set-TextValue "D18" "2.597.15"
Set-TextValue "E18" "  -0.66%  "
Set-TextValue "D19" "362.07"
Set-TextValue "E19" "  -1.28%  "
Set-TextValue "E20" "  -1.36%  "
Set-TextValue "E21" "  -4.28%  "
Set-TextValue "D22" "4.28"
Set-TextValue "E22" "  -0.53%  "
Set-TextValue "D23" "4.82"
Set-TextValue "E23" "  -2.10%  "
Set-TextValue "E24" "  -0.21%  "
Set-TextValue "D25" "72.44"
Set-TextValue "E25" "  +2.97%  "
Set-TextValue "E26" "  +0.05%  "
Set-TextValue "D27" "9.90"
Set-TextValue "E27" "  +0.80%  "
Set-TextValue "E28" "  -0.71%  "
Set-TextValue "B29" "Bittensor"
Set-TextValue "C29" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D29" "581.45"
Set-TextValue "E29" "  +0.83%  "
Set-TextValue "B30" "Binance-PegBSC-USD"
Set-TextValue "C30" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D30" "1.00"
Set-TextValue "E30" "  -0.10%  "
Set-TextValue "E31" "  -5.94%  "
Set-TextValue "E32" "  -4.85%  "
Set-TextValue "E33" "  -3.77%  "
Set-TextValue "E34" "  -3.37%  "
Set-TextValue "D35" "0.999"
Set-TextValue "E35" "  +0.01%  "
Set-TextValue "E36" "  -5.35%  "
Set-TextValue "E37" "  -2.52%  "
Set-TextValue "D38" "156.26"
Set-TextValue "E38" "  -1.22%  "
Set-TextValue "E39" "  -2.26%  "
Set-TextValue "E40" "  -1.30%  "
Set-TextValue "D41" "1.84"
Set-TextValue "E41" "  -0.56%  "
Set-TextValue "D42" "5.17"
Set-TextValue "E42" "  -3.34%  "
Set-TextValue "D43" "17.07"
Set-TextValue "E43" "  +3.92%  "
Set-TextValue "E44" "  -4.29%  "
Set-TextValue "E45" "  -0.13%  "
Set-TextValue "D46" "151.92"
Set-TextValue "E46" "  -2.95%  "
Set-TextValue "E47" "  -0.90%  "
Set-TextValue "E48" "  -1.38%  "
Set-TextValue "D49" "1.67"
Set-TextValue "E49" "  -2.58%  "
Set-TextValue "D50" "0.0776"
Set-TextValue "E50" "  -1.64%  "
Set-TextValue "D51" "21.29"
Set-TextValue "E51" "  +1.74%  "
